# issue #5: stock data output to json file
#
# The 股票 (stock) worksheet gains a new "property_category" column
# (literal value "stock" on every data row), inserted between the
# existing "total" and "date" columns. Everything to the right of the
# insertion point (date / legislator_name / legislator_id) shifts over
# by one column automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a blank column at H (currently "date"), pushing date /
# legislator_name / legislator_id one column to the right (H->I, I->J,
# J->K) and leaving a fresh, empty column H in their place.
$ws.Range("H1:H3").EntireColumn.Insert()

# Header for the freshly inserted column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Populate the new column for every existing data row with "stock".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
